# Daily attendance processing - 2026-01-24 18:43:59
# Reorders the "Recorded By" (column G) entries so that a leading
# "System, " token is moved to become the second item in the
# comma-separated list (e.g. "System, a@b.com" -> "a@b.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row() + $usedRange.Rows().Count() - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value()

    if ($val -ne $null -and $val -is [string] -and $val.StartsWith("System, ")) {
        $parts = @($val -split ", ")
        $rest = @()
        for ($i = 1; $i -lt $parts.Length; $i++) { $rest += $parts[$i] }

        $newParts = @($rest[0], "System")
        for ($i = 1; $i -lt $rest.Length; $i++) { $newParts += $rest[$i] }

        $newVal = $newParts -join ", "
        $cell.Value = $newVal
    }
}
